$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with freshly scraped crypto data.
# Price cells are forced to Text (leading apostrophe) to match the source sheet's
# inline-string storage instead of being auto-parsed as numbers by Excel.

$ws.Range("D2").Value = "'27.971.61"
$ws.Range("E2").Value = "  -3.05%  "
$ws.Range("D3").Value = "'1.861.75"
$ws.Range("E3").Value = "  -2.24%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'318.11"
$ws.Range("E5").Value = "  -1.91%  "
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").Value = "'0.4372"
$ws.Range("E7").Value = "  -4.61%  "
$ws.Range("E8").Value = "  -2.94%  "
$ws.Range("D9").Value = "'0.07512"
$ws.Range("E9").Value = "  -2.55%  "
$ws.Range("D10").Value = "'0.9370"
$ws.Range("E10").Value = "  -4.23%  "
$ws.Range("D11").Value = "'21.30"
$ws.Range("E11").Value = "  -3.77%  "
$ws.Range("D12").Value = "'1.894.00"
$ws.Range("E12").Value = "  +0.44%  "
$ws.Range("D13").Value = "'6.740"
$ws.Range("E13").Value = "  -2.92%  "
$ws.Range("D14").Value = "'5.445"
$ws.Range("E14").Value = "  -3.96%  "
$ws.Range("D15").Value = "'0.06835"
$ws.Range("E15").Value = "  -3.20%  "
$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = "  -0.23%  "
$ws.Range("D17").Value = "'81.63"
$ws.Range("D18").Value = "'0.000009037"
$ws.Range("E18").Value = "  -4.36%  "
$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("D20").Value = "'15.93"
$ws.Range("E20").Value = "  -4.21%  "
$ws.Range("D21").Value = "'27.950.78"
$ws.Range("E21").Value = "  -3.12%  "
$ws.Range("D22").Value = "'5.109"
$ws.Range("E22").Value = "  -3.77%  "
$ws.Range("D23").Value = "'11.06"
$ws.Range("E23").Value = "  +1.27%  "
$ws.Range("D24").Value = "'2.112.56"
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").Value = "'2.000"
$ws.Range("E25").Value = "  -4.71%  "
$ws.Range("D26").Value = "'154.39"
$ws.Range("E26").Value = "  -2.66%  "
$ws.Range("D27").Value = "'18.37"
$ws.Range("E27").Value = "  -3.53%  "
$ws.Range("D28").Value = "'5.430"
$ws.Range("E28").Value = "  -4.06%  "
$ws.Range("D29").Value = "'113.60"
$ws.Range("E29").Value = "  -3.33%  "
$ws.Range("D30").Value = "'1.742"
$ws.Range("E30").Value = "  -6.77%  "
$ws.Range("D31").Value = "'0.08990"
$ws.Range("E31").Value = "  -3.24%  "
$ws.Range("D32").Value = "'0.8116"
$ws.Range("E32").Value = "  -6.00%  "
$ws.Range("D33").Value = "'4.818"
$ws.Range("E33").Value = "  -5.31%  "
$ws.Range("D34").Value = "'1.176"
$ws.Range("E34").Value = "  -5.34%  "
$ws.Range("D35").Value = "'2.942"
$ws.Range("E35").Value = "  -2.81%  "
$ws.Range("E37").Value = "  -3.62%  "
$ws.Range("D38").Value = "'1.117"
$ws.Range("E38").Value = "  -3.58%  "
$ws.Range("D39").Value = "'0.01979"
$ws.Range("E39").Value = "  -3.01%  "
$ws.Range("D40").Value = "'2.895"
$ws.Range("E40").Value = "  -0.29%  "
$ws.Range("D41").Value = "'0.5260"
$ws.Range("E41").Value = "  -4.03%  "
$ws.Range("D42").Value = "'7.052"
$ws.Range("E42").Value = "  -5.32%  "
$ws.Range("D43").Value = "'0.1690"
$ws.Range("E43").Value = "  -3.44%  "
$ws.Range("D44").Value = "'8.808"
$ws.Range("E44").Value = "  -5.57%  "
$ws.Range("D45").Value = "'0.06774"
$ws.Range("E45").Value = "  -1.52%  "
$ws.Range("D46").Value = "'0.4911"
$ws.Range("E46").Value = "  -5.01%  "
$ws.Range("D47").Value = "'10.62"
$ws.Range("E47").Value = "  -5.76%  "
$ws.Range("D48").Value = "'106.38"
$ws.Range("E48").Value = "  -3.61%  "
$ws.Range("E49").Value = "  -5.33%  "
$ws.Range("D50").Value = "'0.9999"
$ws.Range("D51").Value = "'1.901"
$ws.Range("E51").Value = "  -11.91%  "
